$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Cells that switch between text-placeholder and numeric styles ---
# (copy value+format from stable template cells in untouched row 15)
$ws.Range("G15").Copy()
$ws.Range("D14").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E14").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("G14").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("H14").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("J14").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("K14").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("G16").PasteSpecial()
$ws.Range("C15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H16").PasteSpecial()
$ws.Range("E15").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial()
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("D17").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E17").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("C18").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("D18").PasteSpecial()
$ws.Range("C15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("E18").PasteSpecial()
$ws.Range("E15").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("I18").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("C20").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("I20").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("D28").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E28").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("J28").PasteSpecial()
$ws.Range("G15").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("K28").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("L28").PasteSpecial()
$ws.Range("H15").Copy()
$ws.Range("L28").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("G31").PasteSpecial()
$ws.Range("C15").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("H31").PasteSpecial()
$ws.Range("E15").Copy()
$ws.Range("H31").PasteSpecial(-4122)

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 180
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 166.666666666667
$ws.Range("N17").Value = 166.666666666667
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 100
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -80
$ws.Range("M18").Value = -80
$ws.Range("N18").Value = -91.666666666666
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 88.888888888888
$ws.Range("I19").Value = 12
$ws.Range("J19").Value = 6
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = -7.692307692307
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 71.428571428571
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("K20").Value = -66.666666666666
$ws.Range("L20").Value = -87.5
$ws.Range("M20").Value = -80
$ws.Range("N20").Value = -96.969696969697
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 66.666666666666
$ws.Range("F21").Value = 35
$ws.Range("G21").Value = 21
$ws.Range("H21").Value = 66.666666666666
$ws.Range("I21").Value = 22
$ws.Range("J21").Value = 14
$ws.Range("K21").Value = 57.142857142857
$ws.Range("L21").Value = -29.032258064516
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -62.068965517241
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 66.666666666666
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 31
$ws.Range("H24").Value = 32.258064516129
$ws.Range("I24").Value = 31
$ws.Range("J24").Value = 23
$ws.Range("K24").Value = 34.782608695652
$ws.Range("L24").Value = -3.125
$ws.Range("M24").Value = 10.714285714285
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 300
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 63.636363636363
$ws.Range("I25").Value = 13
$ws.Range("J25").Value = 6
$ws.Range("K25").Value = 116.666666666667
$ws.Range("L25").Value = -7.142857142857
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 11
$ws.Range("H26").Value = 45.454545454545
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 100
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 11.111111111111
